# "cleared files/added additional files for automation"
#
# The live user-record rows (2-6) in the batch sheet contained real personal
# data (names, emails, IDs) plus mailto: hyperlinks. For the automation
# template this data needs to be wiped while leaving the header row and the
# existing cell formatting/styling intact - i.e. a "Select the data rows ->
# Clear Contents" pass, plus removal of the now-orphaned hyperlinks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select every data row (row 2 down to the bottom of the sheet) the way a
# user would by dragging across the row headers, then clear their contents.
# This leaves cell styles (s="...") untouched but drops the stored values,
# matching rows 2-6 losing their <v>/t="s" content while keeping s="...".
$dataRows = $ws.Rows("2:1048576")
[void]$dataRows.Select()
$dataRows.ClearContents()

# The cleared cells (D2, A3:A6, D3:D6) carried mailto: hyperlinks to the
# (now gone) user emails - drop those too.
$ws.Hyperlinks.Delete()
